$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: remove the "Skill Highlights" bullet list (heading + its
# seven bullet items, with their spacer paragraphs) from the sidebar.
# The block runs from the blank spacer paragraph right before the
# "Skill Highlights" heading through the "Service-focused" bullet,
# inclusive; the blank spacer paragraph that follows (which precedes
# the "Languages" heading) is left untouched.
# ---------------------------------------------------------------------
$headingIdx = 0
$lastBulletIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text.TrimEnd()
    if ($t -eq "Skill Highlights" -and $p.Range.Font.Bold) {
        $headingIdx = $i
    }
    if ($t -eq "Service-focused") {
        $lastBulletIdx = $i
    }
}

$startPara = $d.Paragraphs($headingIdx - 1)
$afterPara = $d.Paragraphs($lastBulletIdx + 1)
$blockRange = $d.Range($startPara.Range.Start, $afterPara.Range.Start)
$blockRange.Delete()

# ---------------------------------------------------------------------
# Change 2: split "References available on request" so the final word
# is wrapped as its own run flanked by gramStart/gramEnd proofing-error
# markers (as Word's grammar checker does when it flags the phrase).
# ---------------------------------------------------------------------
$refIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text.TrimEnd()
    if ($t -eq "References available on request") {
        $refIdx = $i
    }
}

$refPara = $d.Paragraphs($refIdx)
$paraStart = $refPara.Range.Start
$paraEnd = $refPara.Range.End

$wordRange = $d.Range($paraStart, $paraEnd)
$found = $wordRange.Find.Execute("request", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$xmlFragment = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="FFFFFF" w:themeColor="background1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>request</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$wordRange.InsertXML($xmlFragment)
